# Apply the text corrections described by the commit diff.
$p = $ppt.ActivePresentation

# Slide 4 ("Problemas de otimização"): fix "desde que a possamos" -> "desde que possamos"
$s4 = $p.Slides.Item(4)
$sh4 = $s4.Shapes.Item(2)
$sh4.TextFrame.TextRange.Text = "Em muitos problemas de otimização, o caminho para se atingir o objetivo é irrelevante, desde que possamos conseguir uma solução para o problema em si"

# Slide 16 ("Busca local com pertubação"): fix typo "pertubação" -> "perturbação"
$s16 = $p.Slides.Item(16)
$sh16 = $s16.Shapes.Item(1)
$sh16.TextFrame.TextRange.Runs(2).Text = "perturbação"
